$wb = $excel.ActiveWorkbook
$wsSrc = $wb.Worksheets.Item("Missing")
$wsDst = $wb.Worksheets.Item("WW_props")
$wsSrc.Range("A2:W3").Copy()
$wsDst.Range("A53").PasteSpecial()
Write-Host "Done"
